$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-09 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-10 Thursday", 2) | Out-Null
$d.Content.Find.Execute("94+1=95", $true, $false, $false, $false, $false, $true, 1, $false, "76+2=78", 2) | Out-Null
$d.Content.Find.Execute("17+47=64", $true, $false, $false, $false, $false, $true, 1, $false, "87-30=57", 2) | Out-Null
$d.Content.Find.Execute("58+38=96", $true, $false, $false, $false, $false, $true, 1, $false, "93+4=97", 2) | Out-Null
$d.Content.Find.Execute("48-11=37", $true, $false, $false, $false, $false, $true, 1, $false, "47-7=40", 2) | Out-Null
$d.Content.Find.Execute("37+3=40", $true, $false, $false, $false, $false, $true, 1, $false, "71-20=51", 2) | Out-Null
$d.Content.Find.Execute("3+23=26", $true, $false, $false, $false, $false, $true, 1, $false, "63+17=80", 2) | Out-Null
$d.Content.Find.Execute("14+85=99", $true, $false, $false, $false, $false, $true, 1, $false, "5+31=36", 2) | Out-Null
$d.Content.Find.Execute("67+0=67", $true, $false, $false, $false, $false, $true, 1, $false, "61+30=91", 2) | Out-Null
$d.Content.Find.Execute("23-6=17", $true, $false, $false, $false, $false, $true, 1, $false, "45-5=40", 2) | Out-Null
$d.Content.Find.Execute("34-22=12", $true, $false, $false, $false, $false, $true, 1, $false, "53+37=90", 2) | Out-Null
$d.Content.Find.Execute("59+38=97", $true, $false, $false, $false, $false, $true, 1, $false, "19-5=14", 2) | Out-Null
$d.Content.Find.Execute("93-72=21", $true, $false, $false, $false, $false, $true, 1, $false, "28-9=19", 2) | Out-Null
$d.Content.Find.Execute("25+16=41", $true, $false, $false, $false, $false, $true, 1, $false, "14-8=6", 2) | Out-Null
$d.Content.Find.Execute("6+81=87", $true, $false, $false, $false, $false, $true, 1, $false, "56+33=89", 2) | Out-Null
$d.Content.Find.Execute("68-37=31", $true, $false, $false, $false, $false, $true, 1, $false, "82+2=84", 2) | Out-Null
$d.Content.Find.Execute("3+71=74", $true, $false, $false, $false, $false, $true, 1, $false, "56-48=8", 2) | Out-Null
$d.Content.Find.Execute("7+3=10", $true, $false, $false, $false, $false, $true, 1, $false, "53-0=53", 2) | Out-Null
$d.Content.Find.Execute("7+16=23", $true, $false, $false, $false, $false, $true, 1, $false, "76-41=35", 2) | Out-Null
$d.Content.Find.Execute("90-56=34", $true, $false, $false, $false, $false, $true, 1, $false, "88-18=70", 2) | Out-Null
$d.Content.Find.Execute("72-59=13", $true, $false, $false, $false, $false, $true, 1, $false, "36-22=14", 2) | Out-Null
$d.Content.Find.Execute("97-20=77", $true, $false, $false, $false, $false, $true, 1, $false, "6+43=49", 2) | Out-Null
$d.Content.Find.Execute("38+23=61", $true, $false, $false, $false, $false, $true, 1, $false, "91-49=42", 2) | Out-Null
$d.Content.Find.Execute("84+14=98", $true, $false, $false, $false, $false, $true, 1, $false, "54-43=11", 2) | Out-Null
$d.Content.Find.Execute("60+18=78", $true, $false, $false, $false, $false, $true, 1, $false, "22+16=38", 2) | Out-Null
$d.Content.Find.Execute("47-11=36", $true, $false, $false, $false, $false, $true, 1, $false, "87-74=13", 2) | Out-Null
$d.Content.Find.Execute("98-82=16", $true, $false, $false, $false, $false, $true, 1, $false, "85-18=67", 2) | Out-Null
$d.Content.Find.Execute("47+17=64", $true, $false, $false, $false, $false, $true, 1, $false, "8-4=4", 2) | Out-Null
$d.Content.Find.Execute("75-64=11", $true, $false, $false, $false, $false, $true, 1, $false, "94-22=72", 2) | Out-Null
$d.Content.Find.Execute("94-30=64", $true, $false, $false, $false, $false, $true, 1, $false, "70+2=72", 2) | Out-Null
$d.Content.Find.Execute("35+29=64", $true, $false, $false, $false, $false, $true, 1, $false, "86-66=20", 2) | Out-Null
$d.Content.Find.Execute("13+46=59", $true, $false, $false, $false, $false, $true, 1, $false, "65+24=89", 2) | Out-Null
$d.Content.Find.Execute("97-1=96", $true, $false, $false, $false, $false, $true, 1, $false, "53-4=49", 2) | Out-Null
$d.Content.Find.Execute("9-3=6", $true, $false, $false, $false, $false, $true, 1, $false, "50-17=33", 2) | Out-Null
$d.Content.Find.Execute("71+1=72", $true, $false, $false, $false, $false, $true, 1, $false, "92+4=96", 2) | Out-Null
$d.Content.Find.Execute("6+35=41", $true, $false, $false, $false, $false, $true, 1, $false, "59+22=81", 2) | Out-Null
$d.Content.Find.Execute("4+18=22", $true, $false, $false, $false, $false, $true, 1, $false, "37+22=59", 2) | Out-Null
$d.Content.Find.Execute("12+24=36", $true, $false, $false, $false, $false, $true, 1, $false, "85-69=16", 2) | Out-Null
$d.Content.Find.Execute("43-16=27", $true, $false, $false, $false, $false, $true, 1, $false, "95-13=82", 2) | Out-Null
$d.Content.Find.Execute("39+6=45", $true, $false, $false, $false, $false, $true, 1, $false, "52-42=10", 2) | Out-Null
$d.Content.Find.Execute("49+33=82", $true, $false, $false, $false, $false, $true, 1, $false, "90-17=73", 2) | Out-Null
$d.Content.Find.Execute("74-26=48", $true, $false, $false, $false, $false, $true, 1, $false, "46+24=70", 2) | Out-Null
$d.Content.Find.Execute("12+10=22", $true, $false, $false, $false, $false, $true, 1, $false, "57+31=88", 2) | Out-Null
$d.Content.Find.Execute("93-3=90", $true, $false, $false, $false, $false, $true, 1, $false, "76-48=28", 2) | Out-Null
$d.Content.Find.Execute("20+7=27", $true, $false, $false, $false, $false, $true, 1, $false, "2+84=86", 2) | Out-Null
$d.Content.Find.Execute("37-28=9", $true, $false, $false, $false, $false, $true, 1, $false, "93-84=9", 2) | Out-Null
$d.Content.Find.Execute("7+59=66", $true, $false, $false, $false, $false, $true, 1, $false, "5+21=26", 2) | Out-Null
$d.Content.Find.Execute("0+48=48", $true, $false, $false, $false, $false, $true, 1, $false, "16+75=91", 2) | Out-Null
$d.Content.Find.Execute("29+17=46", $true, $false, $false, $false, $false, $true, 1, $false, "33-27=6", 2) | Out-Null
$d.Content.Find.Execute("79-8=71", $true, $false, $false, $false, $false, $true, 1, $false, "63-6=57", 2) | Out-Null
$d.Content.Find.Execute("0+94=94", $true, $false, $false, $false, $false, $true, 1, $false, "18+9=27", 2) | Out-Null
$d.Content.Find.Execute("57-2=55", $true, $false, $false, $false, $false, $true, 1, $false, "98-92=6", 2) | Out-Null
$d.Content.Find.Execute("37-2=35", $true, $false, $false, $false, $false, $true, 1, $false, "26+46=72", 2) | Out-Null
$d.Content.Find.Execute("46+47=93", $true, $false, $false, $false, $false, $true, 1, $false, "83-42=41", 2) | Out-Null
$d.Content.Find.Execute("11+28=39", $true, $false, $false, $false, $false, $true, 1, $false, "41-19=22", 2) | Out-Null
$d.Content.Find.Execute("55-37=18", $true, $false, $false, $false, $false, $true, 1, $false, "44-40=4", 2) | Out-Null
$d.Content.Find.Execute("62-38=24", $true, $false, $false, $false, $false, $true, 1, $false, "7+29=36", 2) | Out-Null
$d.Content.Find.Execute("40-23=17", $true, $false, $false, $false, $false, $true, 1, $false, "58+14=72", 2) | Out-Null
$d.Content.Find.Execute("14+24=38", $true, $false, $false, $false, $false, $true, 1, $false, "64-49=15", 2) | Out-Null
$d.Content.Find.Execute("17+28=45", $true, $false, $false, $false, $false, $true, 1, $false, "88-58=30", 2) | Out-Null
$d.Content.Find.Execute("96-77=19", $true, $false, $false, $false, $false, $true, 1, $false, "7+76=83", 2) | Out-Null
$d.Content.Find.Execute("92-84=8", $true, $false, $false, $false, $false, $true, 1, $false, "49+6=55", 2) | Out-Null
$d.Content.Find.Execute("88-82=6", $true, $false, $false, $false, $false, $true, 1, $false, "76-66=10", 2) | Out-Null
$d.Content.Find.Execute("89-47=42", $true, $false, $false, $false, $false, $true, 1, $false, "28+56=84", 2) | Out-Null
$d.Content.Find.Execute("10+51=61", $true, $false, $false, $false, $false, $true, 1, $false, "48-48=0", 2) | Out-Null
$d.Content.Find.Execute("62+35=97", $true, $false, $false, $false, $false, $true, 1, $false, "10-9=1", 2) | Out-Null
$d.Content.Find.Execute("23+12=35", $true, $false, $false, $false, $false, $true, 1, $false, "54+19=73", 2) | Out-Null
$d.Content.Find.Execute("40-38=2", $true, $false, $false, $false, $false, $true, 1, $false, "68+1=69", 2) | Out-Null
$d.Content.Find.Execute("63+3=66", $true, $false, $false, $false, $false, $true, 1, $false, "8+66=74", 2) | Out-Null
$d.Content.Find.Execute("54-26=28", $true, $false, $false, $false, $false, $true, 1, $false, "46-22=24", 2) | Out-Null
$d.Content.Find.Execute("8+32=40", $true, $false, $false, $false, $false, $true, 1, $false, "99-95=4", 2) | Out-Null
$d.Content.Find.Execute("83-29=54", $true, $false, $false, $false, $false, $true, 1, $false, "54-3=51", 2) | Out-Null
$d.Content.Find.Execute("15+5=20", $true, $false, $false, $false, $false, $true, 1, $false, "26+3=29", 2) | Out-Null
$d.Content.Find.Execute("10+1=11", $true, $false, $false, $false, $false, $true, 1, $false, "79-61=18", 2) | Out-Null
$d.Content.Find.Execute("8+17=25", $true, $false, $false, $false, $false, $true, 1, $false, "50+40=90", 2) | Out-Null
$d.Content.Find.Execute("45+20=65", $true, $false, $false, $false, $false, $true, 1, $false, "58+13=71", 2) | Out-Null
$d.Content.Find.Execute("30+5=35", $true, $false, $false, $false, $false, $true, 1, $false, "47-19=28", 2) | Out-Null
$d.Content.Find.Execute("84-38=46", $true, $false, $false, $false, $false, $true, 1, $false, "67+4=71", 2) | Out-Null
$d.Content.Find.Execute("85-72=13", $true, $false, $false, $false, $false, $true, 1, $false, "6+38=44", 2) | Out-Null
$d.Content.Find.Execute("48-36=12", $true, $false, $false, $false, $false, $true, 1, $false, "24+5=29", 2) | Out-Null
$d.Content.Find.Execute("34-6=28", $true, $false, $false, $false, $false, $true, 1, $false, "35-25=10", 2) | Out-Null
$d.Content.Find.Execute("83-18=65", $true, $false, $false, $false, $false, $true, 1, $false, "79-4=75", 2) | Out-Null
$d.Content.Find.Execute("67+23=90", $true, $false, $false, $false, $false, $true, 1, $false, "25+67=92", 2) | Out-Null
$d.Content.Find.Execute("24-14=10", $true, $false, $false, $false, $false, $true, 1, $false, "15+58=73", 2) | Out-Null
$d.Content.Find.Execute("32+6=38", $true, $false, $false, $false, $false, $true, 1, $false, "25-17=8", 2) | Out-Null
$d.Content.Find.Execute("95-49=46", $true, $false, $false, $false, $false, $true, 1, $false, "22+35=57", 2) | Out-Null
$d.Content.Find.Execute("56-12=44", $true, $false, $false, $false, $false, $true, 1, $false, "62-6=56", 2) | Out-Null
$d.Content.Find.Execute("28+10=38", $true, $false, $false, $false, $false, $true, 1, $false, "78+2=80", 2) | Out-Null
$d.Content.Find.Execute("39+14=53", $true, $false, $false, $false, $false, $true, 1, $false, "66-60=6", 2) | Out-Null
$d.Content.Find.Execute("75-61=14", $true, $false, $false, $false, $false, $true, 1, $false, "16+61=77", 2) | Out-Null
$d.Content.Find.Execute("73-56=17", $true, $false, $false, $false, $false, $true, 1, $false, "30-29=1", 2) | Out-Null
$d.Content.Find.Execute("57-21=36", $true, $false, $false, $false, $false, $true, 1, $false, "56-37=19", 2) | Out-Null
$d.Content.Find.Execute("13+68=81", $true, $false, $false, $false, $false, $true, 1, $false, "54+13=67", 2) | Out-Null
$d.Content.Find.Execute("6+59=65", $true, $false, $false, $false, $false, $true, 1, $false, "5+89=94", 2) | Out-Null
$d.Content.Find.Execute("3+48=51", $true, $false, $false, $false, $false, $true, 1, $false, "96-55=41", 2) | Out-Null
$d.Content.Find.Execute("92-67=25", $true, $false, $false, $false, $false, $true, 1, $false, "88-48=40", 2) | Out-Null
$d.Content.Find.Execute("18+52=70", $true, $false, $false, $false, $false, $true, 1, $false, "9+28=37", 2) | Out-Null
$d.Content.Find.Execute("69-26=43", $true, $false, $false, $false, $false, $true, 1, $false, "36-16=20", 2) | Out-Null
$d.Content.Find.Execute("8+51=59", $true, $false, $false, $false, $false, $true, 1, $false, "19-10=9", 2) | Out-Null
$d.Content.Find.Execute("1+70=71", $true, $false, $false, $false, $false, $true, 1, $false, "16+46=62", 2) | Out-Null
$d.Content.Find.Execute("21+73=94", $true, $false, $false, $false, $false, $true, 1, $false, "46+5=51", 2) | Out-Null

Write-Output "Replacements applied"
